# Update "Lũy kế tháng CẦN THƠ" report data to reflect latest sync from Notion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-07-24T16:01:00.000Z"

# Row 4
$ws.Range("D4").Value = $newTimestamp
$ws.Range("T4").Value = 82800000
$ws.Range("AA4").Value = 197821000
$ws.Range("AE4").Value = 286350000
$ws.Range("AH4").Value = 250850000

# Row 5
$ws.Range("D5").Value = $newTimestamp
$ws.Range("T5").Value = 65500000
$ws.Range("AA5").Value = 222767000
$ws.Range("AE5").Value = 231700000
$ws.Range("AH5").Value = 186500000

# Row 6
$ws.Range("D6").Value = $newTimestamp
$ws.Range("T6").Value = 26000000
$ws.Range("AA6").Value = 165774000
$ws.Range("AE6").Value = 257900000
$ws.Range("AH6").Value = 242900000

# Row 8
$ws.Range("D8").Value = $newTimestamp
$ws.Range("T8").Value = 35500000
$ws.Range("AA8").Value = 487327000
$ws.Range("AE8").Value = 651800000
$ws.Range("AH8").Value = 582600000

# Row 12
$ws.Range("D12").Value = $newTimestamp
$ws.Range("T12").Value = 130500000
$ws.Range("AA12").Value = 145954000
$ws.Range("AE12").Value = 236900000
$ws.Range("AH12").Value = 217900000

# Row 13
$ws.Range("D13").Value = $newTimestamp
$ws.Range("T13").Value = 35800000
$ws.Range("W13").Value = 88012000
$ws.Range("AA13").Value = 184918000
$ws.Range("AE13").Value = 272930000
$ws.Range("AH13").Value = 227730000
$ws.Range("AK13").Value = 37
$ws.Range("AQ13").Value = 263530000
